# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the per-locale report sheets (zh-cn and de-de), as part of
# regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-11 18:45:45"
$zhcn.Range("H2").Value = "2016-03-11 18:46:04"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-11 18:45:49"
$dede.Range("H2").Value = "2016-03-11 18:46:14"
